$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the three rows whose "Sending cluster" is "Resolving-Mac" (old rows 8-10),
# which also drops the now-unused "Resolving-Mac" shared string.
$ws.Rows("8:10").Delete()

$ws.Range("E2").Value2 = 3
$ws.Range("F2").Value2 = 1
$ws.Range("G2").Value2 = 1.121263666666666
$ws.Range("H2").Value2 = 3.363791
$ws.Range("I2").Value2 = 0.8978163344397481
$ws.Range("J2").Value2 = 0.9294754023256565
$ws.Range("K2").Value2 = 3
$ws.Range("L2").Value2 = 1
$ws.Range("M2").Value2 = 1.219350333333333
$ws.Range("N2").Value2 = 3.658051
$ws.Range("O2").Value2 = 0.2422674834150417
$ws.Range("P2").Value2 = 0.2576244469655636
$ws.Range("Q2").Value2 = 1.367213225704555
$ws.Range("R2").Value2 = 12.304919031341
$ws.Range("S2").Value2 = 0.2175117039136352
$ws.Range("T2").Value2 = 0.239455586492242
$ws.Range("E3").Value2 = 3
$ws.Range("F3").Value2 = 1
$ws.Range("G3").Value2 = 1.121263666666666
$ws.Range("H3").Value2 = 3.363791
$ws.Range("I3").Value2 = 0.8978163344397481
$ws.Range("J3").Value2 = 0.9294754023256565
$ws.Range("K3").Value2 = 3
$ws.Range("L3").Value2 = 1
$ws.Range("M3").Value2 = 2.913661333333334
$ws.Range("N3").Value2 = 8.740984000000001
$ws.Range("O3").Value2 = 0.5789028628226193
$ws.Range("P3").Value2 = 0.6155986258624717
$ws.Range("Q3").Value2 = 3.266982590038222
$ws.Range("R3").Value2 = 29.402843310344
$ws.Range("S3").Value2 = 0.5197484462960803
$ws.Range("T3").Value2 = 0.5721837804446422
$ws.Range("E4").Value2 = 3
$ws.Range("F4").Value2 = 1
$ws.Range("G4").Value2 = 1.121263666666666
$ws.Range("H4").Value2 = 3.363791
$ws.Range("I4").Value2 = 0.8978163344397481
$ws.Range("J4").Value2 = 0.9294754023256565
$ws.Range("K4").Value2 = 2
$ws.Range("L4").Value2 = 1
$ws.Range("M4").Value2 = 0.900063
$ws.Range("N4").Value2 = 1.800126
$ws.Range("O4").Value2 = 0.178829653762339
$ws.Range("P4").Value2 = 0.1267769271719646
$ws.Range("Q4").Value2 = 1.009207939611
$ws.Range("R4").Value2 = 6.055247637665999
$ws.Range("S4").Value2 = 0.1605561842300325
$ws.Range("T4").Value2 = 0.1178360353887722
$ws.Range("E5").Value2 = 2
$ws.Range("F5").Value2 = 1
$ws.Range("G5").Value2 = 0.127615
$ws.Range("H5").Value2 = 0.25523
$ws.Range("I5").Value2 = 0.1021836655602519
$ws.Range("J5").Value2 = 0.07052459767434344
$ws.Range("K5").Value2 = 3
$ws.Range("L5").Value2 = 1
$ws.Range("M5").Value2 = 1.219350333333333
$ws.Range("N5").Value2 = 3.658051
$ws.Range("O5").Value2 = 0.2422674834150417
$ws.Range("P5").Value2 = 0.2576244469655636
$ws.Range("Q5").Value2 = 0.1556073927883333
$ws.Range("R5").Value2 = 0.93364435673
$ws.Range("S5").Value2 = 0.02475577950140651
$ws.Range("T5").Value2 = 0.0181688604733216
$ws.Range("E6").Value2 = 2
$ws.Range("F6").Value2 = 1
$ws.Range("G6").Value2 = 0.127615
$ws.Range("H6").Value2 = 0.25523
$ws.Range("I6").Value2 = 0.1021836655602519
$ws.Range("J6").Value2 = 0.07052459767434344
$ws.Range("K6").Value2 = 3
$ws.Range("L6").Value2 = 1
$ws.Range("M6").Value2 = 2.913661333333334
$ws.Range("N6").Value2 = 8.740984000000001
$ws.Range("O6").Value2 = 0.5789028628226193
$ws.Range("P6").Value2 = 0.6155986258624717
$ws.Range("Q6").Value2 = 0.3718268910533334
$ws.Range("R6").Value2 = 2.23096134632
$ws.Range("S6").Value2 = 0.05915441652653893
$ws.Range("T6").Value2 = 0.04341484541782949
$ws.Range("E7").Value2 = 2
$ws.Range("F7").Value2 = 1
$ws.Range("G7").Value2 = 0.127615
$ws.Range("H7").Value2 = 0.25523
$ws.Range("I7").Value2 = 0.1021836655602519
$ws.Range("J7").Value2 = 0.07052459767434344
$ws.Range("K7").Value2 = 2
$ws.Range("L7").Value2 = 1
$ws.Range("M7").Value2 = 0.900063
$ws.Range("N7").Value2 = 1.800126
$ws.Range("O7").Value2 = 0.178829653762339
$ws.Range("P7").Value2 = 0.1267769271719646
$ws.Range("Q7").Value2 = 0.114861539745
$ws.Range("R7").Value2 = 0.45944615898
$ws.Range("S7").Value2 = 0.0182734695323065
$ws.Range("T7").Value2 = 0.00894089178319234
